$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 19:52"

# Reorder countries: Ghana now appears before Estonia and Bosnia y Herzegovina.
# Row 75 becomes Ghana (new data), row 76 becomes Estonia (old row 75 data),
# row 77 becomes Bosnia y Herzegovina (old row 76 data).
$ws.Range("A75").Value = "Ghana"
$ws.Range("B75").Value = 1671
$ws.Range("C75").Value = 121
$ws.Range("D75").Value = 188
$ws.Range("E75").Value = 1467
$ws.Range("F75").Value = 4
$ws.Range("G75").Value = 5
$ws.Range("H75").Value = 16

$ws.Range("A76").Value = "Estonia"
$ws.Range("B76").Value = 1660
$ws.Range("C76").Value = 13
$ws.Range("D76").Value = 240
$ws.Range("E76").Value = 1370
$ws.Range("F76").Value = 9
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 50

$ws.Range("A77").Value = "Bosnia y Herzegovina"
$ws.Range("B77").Value = 1585
$ws.Range("C77").Value = 20
$ws.Range("D77").Value = 682
$ws.Range("E77").Value = 840
$ws.Range("F77").Value = 4
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 63

# Update statistics for Estados Unidos (row 4)
$ws.Range("B4").Value = 1022259
$ws.Range("C4").Value = 11903
$ws.Range("E4").Value = 824470
$ws.Range("G4").Value = 1065
$ws.Range("H4").Value = 57862

# Update statistics for Irlanda (row 24)
$ws.Range("B24").Value = 19877
$ws.Range("C24").Value = 229
$ws.Range("E24").Value = 9485
$ws.Range("G24").Value = 57
$ws.Range("H24").Value = 1159

# Update statistics for Kazajistan (row 61)
$ws.Range("B61").Value = 3027
$ws.Range("C61").Value = 192
$ws.Range("D61").Value = 774
$ws.Range("E61").Value = 2228

# Update statistics for Georgia (row 109)
$ws.Range("D109").Value = 168
$ws.Range("E109").Value = 337

# Update statistics for Isla de Man (row 125)
$ws.Range("B125").Value = 309
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 252
$ws.Range("E125").Value = 36
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 21
